$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.781.16"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "2.266.49"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.531"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").Value = "2.617.45"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").Value = "2.280.44"
$ws.Range("E17").Value = "  +3.95%  "
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").Value = "41.690.17"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.28%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("E30").Value = "  -4.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0744"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.32%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.116"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.104"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "2.019.95"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.78%  "
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.77%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.30%  "
